# 自动更新Excel文件 - decrement the "剩余" (remaining days) counter in column E
# for each data row. When a row's remaining count hits 1, the cycle restarts:
# remaining resets to 10 and the start date (column F) advances by 10 days.
# Row 36 is left untouched (its start date is a malformed value, so it is
# excluded from the daily update, matching the source data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 99
$skipRows = @(36)

for ($r = 2; $r -le $lastRow; $r++) {
    if ($skipRows -contains $r) { continue }

    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $e = $eCell.Value2
    if ($e -eq $null) { continue }

    if ($e -eq 1) {
        $f = $fCell.Value2
        $eCell.Value = 10
        $fCell.Value = $f + 10
    } else {
        $eCell.Value = $e - 1
    }
}
